$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "37.097.94"
Set-TextValue $ws.Range("E2") "  -0.69%  "
Set-TextValue $ws.Range("D3") "2.076.59"
Set-TextValue $ws.Range("E3") "  -1.00%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "253.10"
Set-TextValue $ws.Range("E5") "  +0.93%  "
Set-TextValue $ws.Range("D6") "0.673"
Set-TextValue $ws.Range("E6") "  +2.11%  "
Set-TextValue $ws.Range("D7") "59.26"
Set-TextValue $ws.Range("E7") "  +14.29%  "
Set-TextValue $ws.Range("E8") "  -0.04%  "
Set-TextValue $ws.Range("D9") "0.389"
Set-TextValue $ws.Range("E9") "  +4.24%  "
Set-TextValue $ws.Range("D10") "61.42"
Set-TextValue $ws.Range("E10") "  -0.46%  "
Set-TextValue $ws.Range("D11") "0.0795"
Set-TextValue $ws.Range("E11") "  +6.69%  "
Set-TextValue $ws.Range("E12") "  +2.54%  "
Set-TextValue $ws.Range("D13") "15.93"
Set-TextValue $ws.Range("E13") "  +5.74%  "
Set-TextValue $ws.Range("D14") "2.383.22"
Set-TextValue $ws.Range("E14") "  -0.85%  "
Set-TextValue $ws.Range("D15") "0.816"
Set-TextValue $ws.Range("E15") "  -2.22%  "
Set-TextValue $ws.Range("D16") "5.54"
Set-TextValue $ws.Range("E16") "  +7.94%  "
Set-TextValue $ws.Range("D17") "2.078.36"
Set-TextValue $ws.Range("E17") "  -0.99%  "
Set-TextValue $ws.Range("D18") "37.072.47"
Set-TextValue $ws.Range("E18") "  -0.70%  "
Set-TextValue $ws.Range("D19") "15.50"
Set-TextValue $ws.Range("E19") "  +10.77%  "
Set-TextValue $ws.Range("D20") "74.36"
Set-TextValue $ws.Range("E20") "  +2.85%  "
Set-TextValue $ws.Range("D21") "0.0₃0922"
Set-TextValue $ws.Range("E21") "  +9.74%  "
Set-TextValue $ws.Range("D22") "5.45"
Set-TextValue $ws.Range("E22") "  +4.33%  "
Set-TextValue $ws.Range("D23") "239.35"
Set-TextValue $ws.Range("E23") "  -0.47%  "
Set-TextValue $ws.Range("D24") "1.00"
Set-TextValue $ws.Range("E24") "  -0.10%  "
Set-TextValue $ws.Range("D25") "2.41"
Set-TextValue $ws.Range("E25") "  -2.08%  "
Set-TextValue $ws.Range("E26") "  +14.11%  "
Set-TextValue $ws.Range("D27") "169.72"
Set-TextValue $ws.Range("E27") "  -0.57%  "
Set-TextValue $ws.Range("D28") "9.30"
Set-TextValue $ws.Range("E28") "  +1.06%  "
Set-TextValue $ws.Range("D29") "20.31"
Set-TextValue $ws.Range("E29") "  -1.95%  "
Set-TextValue $ws.Range("E30") "  +2.53%  "
Set-TextValue $ws.Range("D31") "4.83"
Set-TextValue $ws.Range("E31") "  +7.52%  "
Set-TextValue $ws.Range("D32") "1.12"
Set-TextValue $ws.Range("E32") "  +4.40%  "
Set-TextValue $ws.Range("D33") "0.0636"
Set-TextValue $ws.Range("E33") "  +4.18%  "
Set-TextValue $ws.Range("D34") "4.43"
Set-TextValue $ws.Range("E34") "  +8.12%  "
Set-TextValue $ws.Range("D35") "0.0914"
Set-TextValue $ws.Range("E35") "  -0.23%  "
Set-TextValue $ws.Range("E36") "  -0.03%  "
Set-TextValue $ws.Range("D37") "2.30"
Set-TextValue $ws.Range("E37") "  -1.21%  "
Set-TextValue $ws.Range("D38") "1.77"
Set-TextValue $ws.Range("E38") "  -3.81%  "
Set-TextValue $ws.Range("E39") "  +26.60%  "
Set-TextValue $ws.Range("D40") "1.36"
Set-TextValue $ws.Range("E40") "  +1.75%  "
Set-TextValue $ws.Range("D41") "4.82"
Set-TextValue $ws.Range("E41") "  +38.66%  "
Set-TextValue $ws.Range("D42") "0.0226"
Set-TextValue $ws.Range("E42") "  +1.00%  "
Set-TextValue $ws.Range("D43") "17.75"
Set-TextValue $ws.Range("E43") "  -3.98%  "
Set-TextValue $ws.Range("D44") "1.16"
Set-TextValue $ws.Range("E44") "  -0.27%  "
Set-TextValue $ws.Range("D45") "98.87"
Set-TextValue $ws.Range("E45") "  -0.39%  "
Set-TextValue $ws.Range("E46") "  +2.23%  "
Set-TextValue $ws.Range("D47") "4.50"
Set-TextValue $ws.Range("E47") "  +13.48%  "
Set-TextValue $ws.Range("E48") "  +8.47%  "
Set-TextValue $ws.Range("D49") "2.96"
Set-TextValue $ws.Range("E49") "  -1.56%  "
Set-TextValue $ws.Range("D50") "1.305.97"
Set-TextValue $ws.Range("E50") "  -1.02%  "
Set-TextValue $ws.Range("D51") "6.93"
Set-TextValue $ws.Range("E51") "  -0.89%  "
